$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2225
$ws.Range("F3").Value = 251
$ws.Range("F4").Value = 157
$ws.Range("F5").Value = 153
$ws.Range("F6").Value = 290
$ws.Range("F8").Value = 660
$ws.Range("F9").Value = 487
$ws.Range("F10").Value = 591
$ws.Range("F11").Value = 360
$ws.Range("F12").Value = 50
$ws.Range("F13").Value = 335
$ws.Range("F14").Value = 938
$ws.Range("F15").Value = 200
$ws.Range("F16").Value = 123
$ws.Range("F17").Value = 78
$ws.Range("F19").Value = 12
$ws.Range("F20").Value = 194
$ws.Range("F21").Value = 72

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 61
$ws.Range("F3").Value = 28
$ws.Range("F4").Value = 148
$ws.Range("F6").Value = 161
$ws.Range("F8").Value = 2377
$ws.Range("F10").Value = 12
$ws.Range("F13").Value = 25
$ws.Range("F16").Value = 2177

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 302
$ws.Range("F4").Value = 154

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 61
$ws.Range("F4").Value = 28
$ws.Range("F5").Value = 2225
$ws.Range("F6").Value = 302
$ws.Range("F7").Value = 251
$ws.Range("F8").Value = 157
$ws.Range("F9").Value = 153
$ws.Range("F10").Value = 290
$ws.Range("F11").Value = 148
$ws.Range("F14").Value = 161
$ws.Range("F15").Value = 154
$ws.Range("F16").Value = 660
$ws.Range("F17").Value = 487
$ws.Range("F18").Value = 591
$ws.Range("F19").Value = 360
$ws.Range("F20").Value = 50
$ws.Range("F21").Value = 335
$ws.Range("F22").Value = 938
$ws.Range("F24").Value = 2377
$ws.Range("F26").Value = 12
$ws.Range("F29").Value = 25
$ws.Range("F30").Value = 200
$ws.Range("F31").Value = 123
$ws.Range("F32").Value = 78
$ws.Range("F36").Value = 12
$ws.Range("F37").Value = 194
$ws.Range("F38").Value = 72
$ws.Range("F39").Value = 2177

